$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nsgvs")
$ws.Activate()
$ws.Range("A1").Value = "test"
